$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changed from 45182 to 45184 for all existing data rows (2..215)
for ($r = 2; $r -le 215; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 215 gains an explicit row height (ht="15" customHeight="1"), matching every other row
$ws.Rows.Item(215).RowHeight = 15

# New row 216
$ws.Cells.Item(216, 1).Value = "A 43263-2023"
$ws.Cells.Item(216, 2).Value = 45183
$ws.Cells.Item(216, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(216, 3).Value = 45184
$ws.Cells.Item(216, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(216, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(216, 5).Value = "ÖDESHÖG"
$ws.Cells.Item(216, 7).Value = 1.8
$ws.Cells.Item(216, 8).Value = 0
$ws.Cells.Item(216, 9).Value = 0
$ws.Cells.Item(216, 10).Value = 0
$ws.Cells.Item(216, 11).Value = 0
$ws.Cells.Item(216, 12).Value = 0
$ws.Cells.Item(216, 13).Value = 0
$ws.Cells.Item(216, 14).Value = 0
$ws.Cells.Item(216, 15).Value = 0
$ws.Cells.Item(216, 16).Value = 0
$ws.Cells.Item(216, 17).Value = 0
$ws.Cells.Item(216, 18).WrapText = $true
$ws.Rows.Item(216).RowHeight = 15

# New row 217
$ws.Cells.Item(217, 1).Value = "A 43265-2023"
$ws.Cells.Item(217, 2).Value = 45183
$ws.Cells.Item(217, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(217, 3).Value = 45184
$ws.Cells.Item(217, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(217, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(217, 5).Value = "ÖDESHÖG"
$ws.Cells.Item(217, 7).Value = 2.8
$ws.Cells.Item(217, 8).Value = 0
$ws.Cells.Item(217, 9).Value = 0
$ws.Cells.Item(217, 10).Value = 0
$ws.Cells.Item(217, 11).Value = 0
$ws.Cells.Item(217, 12).Value = 0
$ws.Cells.Item(217, 13).Value = 0
$ws.Cells.Item(217, 14).Value = 0
$ws.Cells.Item(217, 15).Value = 0
$ws.Cells.Item(217, 16).Value = 0
$ws.Cells.Item(217, 17).Value = 0
$ws.Cells.Item(217, 18).WrapText = $true
